$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A86").Value = 221
$ws.Range("B86").Value = 138
$ws.Range("C86").Value = 76
$ws.Range("D86").Value = 4
$ws.Range("E86").Value = 3
$ws.Range("F86").Value = 85
$ws.Range("G86").Value = 80
$ws.Range("H86").Value = 14
$ws.Range("I86").Value = 0
